# ETF.xlsx update
#
# 1. Reorder the worksheet tabs: "Equity" moves from the first position to
#    the last position (after "Bond" and "Alternative"), so the new tab
#    order becomes: Bond, Alternative, Equity.
# 2. Append a new holding to the "Equity" sheet:
#    1329.T | (blank) | iShares Core Nikkei 225 ETF | Tokyo | (blank) | 0 | 0.0182 | 16.89

$wb = $excel.ActiveWorkbook

# --- 1. Move the "Equity" worksheet to the end of the workbook -----------
$equitySheet = $wb.Worksheets.Item("Equity")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$equitySheet.Move($null, $lastSheet)

# --- 2. Append the new ETF row to the "Equity" worksheet -----------------
$equity = $wb.Worksheets.Item("Equity")
$newRow = $equity.UsedRange.Rows.Count + 1

$equity.Cells.Item($newRow, 1).Value = "1329.T"
$equity.Cells.Item($newRow, 3).Value = "iShares Core Nikkei 225 ETF"
$equity.Cells.Item($newRow, 4).Value = "Tokyo"
$equity.Cells.Item($newRow, 6).Value = 0
$equity.Cells.Item($newRow, 7).Value = 0.0182
$equity.Cells.Item($newRow, 8).Value = 16.89
